# Generate Report for Handoff
# Updates the localization-status workbook: status moves from
# "Handed back: in sync with en-US" to "Ready for handoff", refreshes
# the handoff timestamps, records version-mismatch error details, and
# resizes a few columns to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$errorA = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec2bbf4505635998afb0cd1e0c8c5e9d8c9ae6e5/e2e/a99ca549-0f60-4d18-ad94-4b495b53a99d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e15bd761f122204168e6b4fec3c6009c189a5b1d/e2e/a99ca549-0f60-4d18-ad94-4b495b53a99d.md."
$errorB = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec2bbf4505635998afb0cd1e0c8c5e9d8c9ae6e5/e2e/e0b4efcb-6692-4c38-9947-31036f592b9c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e15bd761f122204168e6b4fec3c6009c189a5b1d/e2e/e0b4efcb-6692-4c38-9947-31036f592b9c.md."

# ---- Overview sheet ----
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-07 14:44:16"

$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 14:44:16"

$overview.Columns.Item(5).ColumnWidth = 17
$overview.Columns.Item(6).ColumnWidth = 17

# ---- zh-cn sheet ----
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-07 14:43:55"
$zhcn.Range("P2").Value = $errorA

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-07 14:43:55"
$zhcn.Range("P3").Value = $errorB

$zhcn.Columns.Item(3).ColumnWidth = 17
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-07 14:44:16"
$dede.Range("P2").Value = $errorA

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-07 14:44:16"
$dede.Range("P3").Value = $errorB

$dede.Columns.Item(3).ColumnWidth = 17
$dede.Columns.Item(16).ColumnWidth = 39.17
